# Update time tracking for sprint 5
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New sprint-5 log entries (rows 24-26) ---------------------------------
# Copy the full row format (date / time / text / wrap styles) down from the
# last fully-formatted data row (23) before writing the new values, so the
# new rows pick up the same number formats / alignments as the rest of the
# table.
$ws.Range("A23:G23").Copy()
$ws.Range("A24:G24").PasteSpecial(-4122)
$ws.Range("A23:G23").Copy()
$ws.Range("A25:G25").PasteSpecial(-4122)
$ws.Range("A23:G23").Copy()
$ws.Range("A26:G26").PasteSpecial(-4122)

# Row 24
$ws.Range("A24").Value = 43571
$ws.Range("B24").Value = 0.41666666666666669
$ws.Range("C24").Value = 0.54166666666666663
$ws.Range("E24").Value = "Design"
$ws.Range("F24").Value = "Sprint 5"
$ws.Range("G24").Value = "Designing architecture to store layer data and use that data to update the image export file."
$ws.Rows.Item(24).RowHeight = 29.25

# Row 25
$ws.Range("A25").Value = 43573
$ws.Range("B25").Value = 0.51041666666666663
$ws.Range("C25").Value = 0.64722222222222225
$ws.Range("E25").Value = "Testing"
$ws.Range("F25").Value = "Sprint 5"
$ws.Range("G25").Value = "Fix bugs associated with larger brush event handlers. Integrated brush handlers to edit image data."
$ws.Rows.Item(25).RowHeight = 45.75

# Row 26
$ws.Range("A26").Value = 43579
$ws.Range("B26").Value = 0.59027777777777779
$ws.Range("C26").Value = 0.68194444444444446
$ws.Range("E26").Value = "Testing"
$ws.Range("F26").Value = "Sprint 5"
$ws.Range("G26").Value = "Final testing and debug of sprint integration, including event handlers for canvas resize and new canvas methods. "
$ws.Rows.Item(26).RowHeight = 73.5

# --- Sprint 5 rollup (I7) + Project Total row (H9/I9) ----------------------
$ws.Range("I7").Formula = "=SUM(D24:D26)"
$ws.Range("I7").Style = $ws.Range("I6").Style

$ws.Range("H9").Value = "Project Total"
$ws.Range("I9").NumberFormat = $ws.Range("I6").NumberFormat

# --- Column G / H sizing ----------------------------------------------------
# G keeps its existing "wrap text" column style; H grows a little to fit the
# new "Project Total" label.
$ws.Columns.Item(8).ColumnWidth = 12.6

# --- Misc view state ---------------------------------------------------------
$ws.Range("K12").Select()
